# Update the lattice-multiplication exercise table: each of the 15 table
# cells keeps its layout (problem line, spaced-digits line, "----" line,
# two "<digit>|    |" lines) but gets new numbers, per the target diff.
#
# Each cell's visible text is a single run of five "lines" joined by
# manual line breaks (<w:br/> == chr(11) in Range.Text), so the safest
# way to reproduce the exact OOXML shape is to overwrite Cell.Range.Text
# wholesale with the same five-line layout, new digits substituted in.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$vt = [char]11   # manual line break, as it appears in Range.Text

# Row/Col -> new content (top line, spaced-digits line, two product lines).
# The "  ----" separator line is unchanged in every cell.
$cellData = @(
    @{Row=1; Col=1; L1="95 x 55"; L2="  5    5"; L3="9|    |"; L4="5|    |"}
    @{Row=1; Col=2; L1="84 x 13"; L2="  1    3"; L3="8|    |"; L4="4|    |"}
    @{Row=1; Col=3; L1="44 x 26"; L2="  2    6"; L3="4|    |"; L4="4|    |"}
    @{Row=2; Col=1; L1="27 x 99"; L2="  9    9"; L3="2|    |"; L4="7|    |"}
    @{Row=2; Col=2; L1="62 x 30"; L2="  3    0"; L3="6|    |"; L4="2|    |"}
    @{Row=2; Col=3; L1="67 x 84"; L2="  8    4"; L3="6|    |"; L4="7|    |"}
    @{Row=3; Col=1; L1="91 x 84"; L2="  8    4"; L3="9|    |"; L4="1|    |"}
    @{Row=3; Col=2; L1="42 x 19"; L2="  1    9"; L3="4|    |"; L4="2|    |"}
    @{Row=3; Col=3; L1="89 x 87"; L2="  8    7"; L3="8|    |"; L4="9|    |"}
    @{Row=4; Col=1; L1="49 x 68"; L2="  6    8"; L3="4|    |"; L4="9|    |"}
    @{Row=4; Col=2; L1="18 x 38"; L2="  3    8"; L3="1|    |"; L4="8|    |"}
    @{Row=4; Col=3; L1="42 x 77"; L2="  7    7"; L3="4|    |"; L4="2|    |"}
    @{Row=5; Col=1; L1="34 x 94"; L2="  9    4"; L3="3|    |"; L4="4|    |"}
    @{Row=5; Col=2; L1="11 x 80"; L2="  8    0"; L3="1|    |"; L4="1|    |"}
    @{Row=5; Col=3; L1="84 x 43"; L2="  4    3"; L3="8|    |"; L4="4|    |"}
)

foreach ($item in $cellData) {
    $cell = $tbl.Cell($item.Row, $item.Col)
    $newText = $item.L1 + $vt + $item.L2 + $vt + "  ----" + $vt + $item.L3 + $vt + $item.L4
    $cell.Range.Text = $newText
}

Write-Output ("Updated " + $cellData.Count + " cells")
